# Update the value in A2 from "Emre Abale" to "Rob Oudman"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Rob Oudman"

# Update the selected cell/active cell from B4 to A2
$ws.Range("A2").Select()
